# use_case_diagram.pptx edit
#
#  - Shape "椭圆 59" (ellipse, id 60): nudge Top up slightly
#      (710674 -> 710366 EMU) and change its text from "登录" to "登出".
#  - Shape "直接箭头连接符 37" (connector, id 38): nudge Top up slightly
#      (1020917 -> 1020609 EMU).
#  - Shape "直接箭头连接符 41" (connector, id 42): nudge Top up slightly
#      (1020917 -> 1020609 EMU) and increase Height slightly
#      (1608194 -> 1608502 EMU).
#
# PowerPoint's COM Shape.Top/.Height/.Left/.Width are expressed in points
# and stored as single-precision (float32) floats, and the EMU value
# written back into the OOXML is obtained by truncating pt*12700 rather
# than rounding it. A naive EMU/12700 division can therefore land one
# EMU below the intended value once it has been through the float32
# round-trip. ConvertTo-SafePt searches for a point value that survives
# that round-trip and reproduces the exact target EMU, so the written
# offsets/extents match the target file bit-for-bit.
function ConvertTo-SafePt {
    param([double]$TargetEmu)
    $EmuPerPt = 12700.0
    $basePt = $TargetEmu / $EmuPerPt
    for ($i = 0; $i -le 20000; $i++) {
        $cand = $basePt + ($i * 0.0000001)
        $f32 = [single]$cand
        $emu = [int]([double]$f32 * $EmuPerPt)
        if ($emu -eq $TargetEmu) {
            return $cand
        }
    }
    return $basePt
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$ellipse = $s.Shapes.Item("椭圆 59")
$ellipse.Top = ConvertTo-SafePt 710366
$ellipse.TextFrame.TextRange.Text = "登出"

$conn1 = $s.Shapes.Item("直接箭头连接符 37")
$conn1.Top = ConvertTo-SafePt 1020609

$conn2 = $s.Shapes.Item("直接箭头连接符 41")
$conn2.Top = ConvertTo-SafePt 1020609
$conn2.Height = ConvertTo-SafePt 1608502
